$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 716, shifting existing rows 716:757 down to 717:758.
$ws.Rows.Item(716).Insert()

# Fill the new row 716 with the new entry: 2026/01/26 (月), time 12, ranking 201.
# Pre-format column A as Text so the date-like string "2026/01/26" is stored as a
# literal string (matching the rest of the sheet) instead of being auto-converted
# into a date serial number, then clear the format again so the cell ends up with
# no explicit style - same as its sibling cells.
$ws.Cells.Item(716, 1).NumberFormat = "@"
$ws.Cells.Item(716, 1).Value = "2026/01/26"
$ws.Cells.Item(716, 1).ClearFormats()

$ws.Cells.Item(716, 2).Value = "月"
$ws.Cells.Item(716, 3).Value = 12
$ws.Cells.Item(716, 4).Value = 201
